$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new docente was added to the "Docentes responsaveis" list, just above
# the existing "Cristina Bormio Nunes" row. Insert a new row at 13 - this
# shifts the existing rows 13-25 (and their heights/formatting) down to
# 14-26.
$ws.Rows("13:13").Insert()

# The row-insert carries column-A's (bold) formatting into the new row;
# this row has no column-A entry at all, so drop it.
$ws.Range("A13").Clear()

# Fill in the new person's entry in columns B and C (same text is mirrored
# into both, like every other row in this "current / modified" table).
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"

# Re-use the exact formatting of the equivalent cells one row below
# (B: wrapped plain text, C: wrapped red text marking new/changed data)
# instead of leaving the inherited bold header formatting in place.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
